# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" sheets, reflecting the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2074
$ws1.Range("F5").Value = 357
$ws1.Range("F6").Value = 616
$ws1.Range("F9").Value = 10655
$ws1.Range("F12").Value = 284
$ws1.Range("F15").Value = 7522
$ws1.Range("F18").Value = 255

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2074
$ws4.Range("F5").Value = 357
$ws4.Range("F6").Value = 616
$ws4.Range("F12").Value = 10655
$ws4.Range("F15").Value = 284
$ws4.Range("F18").Value = 7522
$ws4.Range("F21").Value = 255
